{"js": "// Add a new sentence (as its own run) right after the existing sentence\n// \"...you can easily switch between relevant charts.\" in the DATASET\n// VISUALS section, describing the new customizable scatter plot.\n\nconst body = context.document.body;\nconst anchorText = \"The visuals from this report can be viewed from this app, and you can easily switch between relevant charts.\";\nconst addition = \" There is also a customizable scatter plot where you can select data based on year and season, and pick attributes as x-axis or y-axis.\";\n\nconst results = body.search(anchorText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor sentence not found\");\n}\n\nconst anchorRange = results.items[0];\n\n// Insert the new sentence right after the matched sentence.\nconst newRange = anchorRange.insertText(addition, Word.InsertLocation.after);\n\n// Give the newly inserted text its own run properties (matches the rest of\n// the paragraph's formatting: Times New Roman, 12pt) so it is written out\n// as a distinct <w:r> rather than being merged into the previous run.\nnewRange.font.name = \"Times New Roman\";\nnewRange.font.nameBidirectional = \"Times New Roman\";\nnewRange.font.size = 12;\n\nawait context.sync();\n", "ps1": "# Add a new sentence (as its own run) right after the existing sentence\n# \"...you can easily switch between relevant charts.\" in the DATASET\n# VISUALS section, describing the new customizable scatter plot.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"The visuals from this report can be viewed from this app, and you can easily switch between relevant charts.\"\n$addition = \" There is also a customizable scatter plot where you can select data based on year and season, and pick attributes as x-axis or y-axis.\"\n\n$rng = $d.Content\n$rng.Find.Execute($anchorText) | Out-Null\n\n# Collapse to the end of the found sentence, then insert the new sentence.\n$rng.Collapse(0) | Out-Null\n$rng.InsertAfter($addition)\n\n# Force the inserted text to carry its own (explicit) run formatting -\n# matching the rest of the paragraph (Times New Roman, 12pt) - so it is\n# written out as a separate <w:r> rather than merged into the prior run.\n$rng.Font.Name = \"Times New Roman\"\n$rng.Font.NameBi = \"Times New Roman\"\n$rng.Font.Size = 12\n"}
